$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Update "Inscritos" (column E) and related counts to reflect new registration data.

# Row 5 (Agropecuaria item): Inscritos 27 -> 29
$ws.Range("E5").Value = 29

# Row 6: Inscritos 56 -> 57
$ws.Range("E6").Value = 57

# Row 12: Pagos 13 -> 14, Inscrições homologadas 15 -> 16
$ws.Range("F12").Value = 14
$ws.Range("H12").Value = 16

# Row 16: Inscritos 319 -> 320
$ws.Range("E16").Value = 320

# Row 17: Inscritos 25 -> 26, Pagos 11 -> 12, Inscrições homologadas 12 -> 13
$ws.Range("E17").Value = 26
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = 13
